$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, pushing existing rows 72..131 down to 73..132.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new price record.
$ws.Cells.Item(72, 1).Value = 11
$ws.Cells.Item(72, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(72, 3).Value = 'Bíobío'
$ws.Cells.Item(72, 4).Value = 44827
$ws.Cells.Item(72, 5).Value = 8
$ws.Cells.Item(72, 6).Value = 'Fruta'
$ws.Cells.Item(72, 7).Value = 100108
$ws.Cells.Item(72, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(72, 9).Value = 100108002
$ws.Cells.Item(72, 10).Value = 'Mango'
$ws.Cells.Item(72, 11).Value = 'Sin especificar'
$ws.Cells.Item(72, 12).Value = 'Primera'
$ws.Cells.Item(72, 13).Value = 180
$ws.Cells.Item(72, 14).Value = 8500
$ws.Cells.Item(72, 15).Value = 9000
$ws.Cells.Item(72, 16).Value = 8778
$ws.Cells.Item(72, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(72, 18).Value = 'Brasil'
$ws.Cells.Item(72, 19).Value = 2194
$ws.Cells.Item(72, 20).Value = 4
